$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = 122
$ws.Range("E4").Value = 79
$ws.Range("D5").Value = "9.6 (5.5)"
$ws.Range("E5").Value = "7.5 (6.4)"
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "0.021"
$ws.Range("F5").Style = "Normal"
$ws.Range("D6").Value = "65 (53.7)"
$ws.Range("E6").Value = "30 (38.5)"
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "0.050"
$ws.Range("F6").Style = "Normal"
$ws.Range("D7").Value = "56 (46.3)"
$ws.Range("E7").Value = "48 (61.5)"
$ws.Range("D8").Value = "48 (39.3)"
$ws.Range("E8").Value = "39 (49.4)"
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "0.209"
$ws.Range("F8").Style = "Normal"
$ws.Range("D9").Value = "74 (60.7)"
$ws.Range("E9").Value = "40 (50.6)"
$ws.Range("D10").Value = "87 (71.9)"
$ws.Range("E10").Value = "56 (71.8)"
$ws.Range("F10").NumberFormat = "@"
$ws.Range("F10").Value = "0.721"
$ws.Range("F10").Style = "Normal"
$ws.Range("D11").Value = "21 (17.4)"
$ws.Range("E11").Value = "11 (14.1)"
$ws.Range("E12").Value = "1 (1.3)"
$ws.Range("D13").Value = "1 (0.8)"
$ws.Range("E13").Value = "1 (1.3)"
$ws.Range("D14").Value = "12 (9.9)"
$ws.Range("E14").Value = "9 (11.5)"
$ws.Range("D15").Value = "11 (9.1)"
$ws.Range("E15").Value = "14 (17.9)"
$ws.Range("F15").NumberFormat = "@"
$ws.Range("F15").Value = "0.105"
$ws.Range("F15").Style = "Normal"
$ws.Range("D16").Value = "110 (90.9)"
$ws.Range("E16").Value = "64 (82.1)"
$ws.Range("D17").Value = "43 (36.8)"
$ws.Range("E17").Value = "33 (45.8)"
$ws.Range("F17").NumberFormat = "@"
$ws.Range("F17").Value = "0.278"
$ws.Range("F17").Style = "Normal"
$ws.Range("D18").Value = "74 (63.2)"
$ws.Range("E18").Value = "39 (54.2)"
$ws.Range("D19").Value = "53 (43.4)"
$ws.Range("E19").Value = "35 (44.9)"
$ws.Range("F19").NumberFormat = "@"
$ws.Range("F19").Value = "0.958"
$ws.Range("F19").Style = "Normal"
$ws.Range("D20").Value = "69 (56.6)"
$ws.Range("E20").Value = "43 (55.1)"
$ws.Range("D21").Value = "55.5 (25.0)"
$ws.Range("E21").Value = "67.3 (25.2)"
$ws.Range("F21").NumberFormat = "@"
$ws.Range("F21").Value = "0.003"
$ws.Range("F21").Style = "Normal"
$ws.Range("D22").Value = "30 (24.6)"
$ws.Range("E22").Value = "21 (26.6)"
$ws.Range("D23").Value = "40 (32.8)"
$ws.Range("E23").Value = "47 (59.5)"
$ws.Range("D24").Value = "52 (42.6)"
$ws.Range("E24").Value = "11 (13.9)"
$ws.Range("D25").Value = "97 (79.5)"
$ws.Range("E25").Value = "62 (78.5)"
$ws.Range("D26").Value = "25 (20.5)"
$ws.Range("E26").Value = "17 (21.5)"
$ws.Range("D27").Value = "21 (17.4)"
$ws.Range("E27").Value = "10 (12.8)"
$ws.Range("F27").NumberFormat = "@"
$ws.Range("F27").Value = "0.509"
$ws.Range("F27").Style = "Normal"
$ws.Range("D28").Value = "100 (82.6)"
$ws.Range("E28").Value = "68 (87.2)"
